$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "分析一下你的哪些性格特点使得你在刚才的那些负面回忆中感到如此的不开心?"
$ws.Range("B3").Value = "刚才的那些负面回忆反映出你是怎样的一个人？你和别人之间有哪些相同点和不同点？"
$ws.Range("B4").Value = "为什么刚才的那些负面回忆偏偏发生在我的身上，而不是别人？"
$ws.Range("B5").Value = "在刚才的那些负面回忆中，我为什么不能把事情处理得更好？"

$ws.Range("D7").Select()
